$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 19:48:00"
$ws.Range("O2").Value = "0.2 °C"
$ws.Range("E3").Value = "2026-02-06 19:48:03"
$ws.Range("H3").Value = "'70%"
$ws.Range("O3").Value = "-2.1 °C"
$ws.Range("E4").Value = "2026-02-06 19:48:06"
$ws.Range("H4").Value = "'56%"
$ws.Range("J4").Value = "997.5 hPa"
$ws.Range("E5").Value = "2026-02-06 19:48:08"
$ws.Range("J5").Value = "997.7 hPa"
$ws.Range("E6").Value = "2026-02-06 19:48:11"
$ws.Range("J6").Value = "998.9 hPa"
$ws.Range("E7").Value = "2026-02-06 19:48:14"
$ws.Range("J7").Value = "998.5 hPa"
$ws.Range("E8").Value = "2026-02-06 19:48:16"
$ws.Range("O8").Value = "10.1 °C"
$ws.Range("E9").Value = "2026-02-06 19:48:19"
$ws.Range("H9").Value = "'85%"
$ws.Range("E10").Value = "2026-02-06 19:48:22"
$ws.Range("E11").Value = "2026-02-06 19:48:24"
$ws.Range("J11").Value = "999.0 hPa"
$ws.Range("O11").Value = "5.3 °C"
$ws.Range("E12").Value = "2026-02-06 19:48:27"
$ws.Range("N12").Value = "7.3 °C 19:14 TU"
$ws.Range("O12").Value = "13.8 °C"
$ws.Range("E13").Value = "2026-02-06 19:48:29"
$ws.Range("O13").Value = "10.3 °C"
$ws.Range("E14").Value = "2026-02-06 19:48:32"
$ws.Range("E15").Value = "2026-02-06 19:48:35"
$ws.Range("J15").Value = "997.9 hPa"
$ws.Range("O15").Value = "10.7 °C"
$ws.Range("E16").Value = "2026-02-06 19:48:37"
$ws.Range("E17").Value = "2026-02-06 19:48:40"
$ws.Range("J17").Value = "999.1 hPa"
$ws.Range("E18").Value = "2026-02-06 19:48:43"
$ws.Range("N18").Value = "-6.6 °C 19:16 TU"
$ws.Range("O18").Value = "-4.6 °C"
$ws.Range("E19").Value = "2026-02-06 19:48:45"
$ws.Range("J19").Value = "1000.0 hPa"
$ws.Range("O19").Value = "9.9 °C"
$ws.Range("E20").Value = "2026-02-06 19:48:48"
$ws.Range("H20").Value = "'80%"
$ws.Range("E21").Value = "2026-02-06 19:48:50"
$ws.Range("J21").Value = "998.1 hPa"
$ws.Range("E22").Value = "2026-02-06 19:48:53"
$ws.Range("O22").Value = "10.6 °C"
$ws.Range("E23").Value = "2026-02-06 19:48:56"
$ws.Range("H23").Value = "'83%"
$ws.Range("J23").Value = "997.9 hPa"
$ws.Range("O23").Value = "10.1 °C"
$ws.Range("E24").Value = "2026-02-06 19:48:59"
$ws.Range("J24").Value = "997.3 hPa"
$ws.Range("E25").Value = "2026-02-06 19:49:01"
$ws.Range("J25").Value = "998.6 hPa"
$ws.Range("E26").Value = "2026-02-06 19:49:04"
$ws.Range("E27").Value = "2026-02-06 19:49:07"
$ws.Range("J27").Value = "997.9 hPa"
$ws.Range("E28").Value = "2026-02-06 19:49:09"
$ws.Range("J28").Value = "1000.0 hPa"
$ws.Range("E29").Value = "2026-02-06 19:49:12"
$ws.Range("K29").Value = "12.1 MJ/m2"
$ws.Range("E30").Value = "2026-02-06 19:49:15"
$ws.Range("H30").Value = "'76%"
$ws.Range("E31").Value = "2026-02-06 19:49:17"
$ws.Range("H31").Value = "'83%"
$ws.Range("J31").Value = "999.3 hPa"
$ws.Range("E32").Value = "2026-02-06 19:49:20"
$ws.Range("O32").Value = "15.5 °C"
$ws.Range("E33").Value = "2026-02-06 19:49:22"
$ws.Range("H33").Value = "'85%"
$ws.Range("E34").Value = "2026-02-06 19:49:25"
$ws.Range("O34").Value = "9.0 °C"
$ws.Range("E35").Value = "2026-02-06 19:49:28"
$ws.Range("G35").Value = "202 cm"
$ws.Range("I35").Value = "0.4 mm"
$ws.Range("O35").Value = "-2.1 °C"
$ws.Range("E36").Value = "2026-02-06 19:49:30"
$ws.Range("I36").Value = "1.7 mm"
$ws.Range("J36").Value = "1000.1 hPa"
$ws.Range("O36").Value = "12.9 °C"
